# Add team record (Wins/Losses/Ties) columns to the data sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells go right after the existing last column (AC) -> AD, AE, AF.
# Copy the existing header style (bold, bordered, centered) from A1 so the
# new headers match the look of the rest of row 1.
$ws.Range("A1").Copy($ws.Range("AD1"))
$ws.Range("A1").Copy($ws.Range("AE1"))
$ws.Range("A1").Copy($ws.Range("AF1"))

$ws.Range("AD1").Value2 = "Wins"
$ws.Range("AE1").Value2 = "Losses"
$ws.Range("AF1").Value2 = "Ties"

# Every player row (2-44) gets the same team record numbers.
for ($r = 2; $r -le 44; $r++) {
    $ws.Cells.Item($r, 30).Value2 = 71   # AD -> Wins
    $ws.Cells.Item($r, 31).Value2 = 91   # AE -> Losses
    $ws.Cells.Item($r, 32).Value2 = 0    # AF -> Ties
}
